$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric "Dirección" values for the existing placeholder rows ---
$ws.Range("A4").Value = 6
$ws.Range("A5").Value = 16
$ws.Range("A6").Value = 26
$ws.Range("A7").Value = 36
$ws.Range("A8").Value = 40
$ws.Range("A9").Value = 44

# --- "Bloque" short codes (column B), written in the same order the
#     original author typed them: FV/HV rows first, then D1/D2/D3 rows ---
$ws.Range("B7").Value = "FV"
$ws.Range("B8").Value = "HV"
$ws.Range("B4").Value = "D1"
$ws.Range("B5").Value = "D2"
$ws.Range("B6").Value = "D3"

# --- "Tamaño (bytes)" values (column C) ---
$ws.Range("C4").Value = 10
$ws.Range("C5").Value = 10
$ws.Range("C6").Value = 10
$ws.Range("C7").Value = 4
$ws.Range("C8").Value = 4

# --- "Descripción" values (column D), D1/D2/D3 rows first, then FV/HV
#     rows, and finally the two descriptions for the pre-existing rows ---
$ws.Range("D4").Value = "Destinatario 1"
$ws.Range("D5").Value = "Destinatario 2"
$ws.Range("D6").Value = "Destinatario 3"
$ws.Range("D7").Value = "Firmware Versión"
$ws.Range("D8").Value = "Hardware Versión"
$ws.Range("D2").Value = "Si esta activado o no"
$ws.Range("D3").Value = "devID del dispositivo"

# --- Row 10 no longer holds any data; clear it so the sheet shrinks to A1:D9 ---
$ws.Cells.Item(10, 1).ClearContents() | Out-Null

# --- Column widths: columns B and D are now best-fit to their (wider) contents ---
$ws.Columns.Item(2).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(4).EntireColumn.AutoFit() | Out-Null

# --- Update the selected cell shown when the sheet is re-opened ---
$ws.Range("E12").Select() | Out-Null
